$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the "last updated" date (C1) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- "FPIEBP" sheet: re-prioritize hard coal (row 3) production/imports/exports ---
$wsFpiebp = $wb.Worksheets.Item("FPIEBP")
$wsFpiebp.Range("B3").Value = 1
$wsFpiebp.Range("C3").Value = 3
$wsFpiebp.Range("D3").Value = 2

# Keep FPIEBP the active/visible tab, with the cursor now resting on E3
$wsFpiebp.Activate() | Out-Null
$wsFpiebp.Range("E3").Select() | Out-Null
